$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 4 ---
# D4: 28 -> 40
$ws.Range("D4").Value = 40
# E4: "105" -> "105|104|103|101"
$ws.Range("E4").Value = "105|104|103|101"

# --- Row 5: blank spacer row, styled like row 3 ---
$ws.Range("A5:F5").Value = "x"
$ws.Range("A5:F5").ClearContents()

# --- Row 6: new user record, styled like row 4 ---
$ws.Range("A4:F4").Copy() | Out-Null
$ws.Range("A6:F6").PasteSpecial() | Out-Null

$ws.Range("A6").Value = "GGG"
$ws.Range("B6").Value = "GGG"
$ws.Range("C6").Value = "835DF613"
$ws.Range("D6").Value = 18
$ws.Range("E6").Value = "104|101"
# F6 already copied as "TRUE" from F4, matching the target value

$excel.CutCopyMode = $false
